$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "30.324.30"
$ws.Range("E2").Value = "  -0.06%  "
$ws.Range("D3").Value = "1.933.62"
$ws.Range("E3").Value = "  +0.08%  "
$ws.Range("D4").Value = "'0.9997"
$ws.Range("E4").Value = "  -0.13%  "
$ws.Range("D5").Value = "'0.7590"
$ws.Range("E5").Value = "  +5.67%  "
$ws.Range("D6").Value = "'244.84"
$ws.Range("E6").Value = "  -2.65%  "
$ws.Range("D7").Value = "'0.9991"
$ws.Range("E7").Value = "  -0.19%  "
$ws.Range("D8").Value = "'0.3184"
$ws.Range("E8").Value = "  -2.51%  "
$ws.Range("E9").Value = "  +0.38%  "
$ws.Range("D10").Value = "'0.06997"
$ws.Range("E10").Value = "  -2.67%  "
$ws.Range("D11").Value = "'0.7783"
$ws.Range("E11").Value = "  -2.61%  "
$ws.Range("E12").Value = "  -1.16%  "
$ws.Range("D13").Value = "1.933.23"
$ws.Range("E13").Value = "  +0.08%  "
$ws.Range("D14").Value = "'5.350"
$ws.Range("E14").Value = "  -1.36%  "
$ws.Range("D15").Value = "'94.32"
$ws.Range("E15").Value = "  -0.38%  "
$ws.Range("D16").Value = "'14.41"
$ws.Range("E16").Value = "  -2.42%  "
$ws.Range("D17").Value = "30.331.39"
$ws.Range("E17").Value = "  +0.02%  "
$ws.Range("D18").Value = "'252.84"
$ws.Range("E18").Value = "  +0.32%  "
$ws.Range("D19").Value = "'0.000007928"
$ws.Range("E19").Value = "  -2.52%  "
$ws.Range("D20").Value = "'5.745"
$ws.Range("E20").Value = "  -0.80%  "
$ws.Range("D21").Value = "2.187.79"
$ws.Range("E21").Value = "  +0.32%  "
$ws.Range("D22").Value = "'0.9986"
$ws.Range("E22").Value = "  -0.24%  "
$ws.Range("D23").Value = "'0.9998"
$ws.Range("E23").Value = "  -0.08%  "
$ws.Range("D24").Value = "'6.680"
$ws.Range("E24").Value = "  -3.48%  "
$ws.Range("D25").Value = "'9.492"
$ws.Range("E25").Value = "  -2.37%  "
$ws.Range("D26").Value = "'165.75"
$ws.Range("E26").Value = "  +0.31%  "
$ws.Range("D27").Value = "'0.1336"
$ws.Range("E27").Value = "  +4.44%  "
$ws.Range("D28").Value = "'18.98"
$ws.Range("E28").Value = "  -1.27%  "
$ws.Range("D29").Value = "'2.182"
$ws.Range("E29").Value = "  -5.99%  "
$ws.Range("D30").Value = "'1.366"
$ws.Range("E30").Value = "  +0.42%  "
$ws.Range("D31").Value = "'1.512"
$ws.Range("E31").Value = "  -2.15%  "
$ws.Range("D32").Value = "'4.382"
$ws.Range("E32").Value = "  -1.02%  "
$ws.Range("D33").Value = "'4.124"
$ws.Range("E33").Value = "  -1.74%  "
$ws.Range("D34").Value = "'0.05162"
$ws.Range("E34").Value = "  -0.93%  "
$ws.Range("D35").Value = "'1.286"
$ws.Range("E35").Value = "  +1.45%  "
$ws.Range("D36").Value = "'0.7498"
$ws.Range("E36").Value = "  +0.36%  "
$ws.Range("D37").Value = "'2.768"
$ws.Range("E37").Value = "  -0.05%  "
$ws.Range("E38").Value = "  -0.17%  "
$ws.Range("D39").Value = "'2.804"
$ws.Range("E39").Value = "  +0.15%  "
$ws.Range("D40").Value = "'77.55"
$ws.Range("E40").Value = "  -1.86%  "
$ws.Range("D41").Value = "'6.415"
$ws.Range("E41").Value = "  -0.40%  "
$ws.Range("D42").Value = "'0.4461"
$ws.Range("E42").Value = "  -1.41%  "
$ws.Range("D43").Value = "'1.965"
$ws.Range("E43").Value = "  -3.12%  "
$ws.Range("E44").Value = "  -0.17%  "
$ws.Range("D45").Value = "'0.8331"
$ws.Range("E45").Value = "  -1.01%  "
$ws.Range("D46").Value = "'100.69"
$ws.Range("E46").Value = "  -1.03%  "
$ws.Range("D47").Value = "'9.750"
$ws.Range("E47").Value = "  -0.23%  "
$ws.Range("D48").Value = "'7.466"
$ws.Range("E48").Value = "  +0.78%  "
$ws.Range("D49").Value = "'982.59"
$ws.Range("E49").Value = "  +11.25%  "
$ws.Range("D50").Value = "'37.37"
$ws.Range("E50").Value = "  +2.29%  "
$ws.Range("D51").Value = "'0.06007"
$ws.Range("E51").Value = "  -1.12%  "
